$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 15
$ws.Range("E15").Value = 103

# Row 38
$ws.Range("E38").Value = 38

# Row 40
$ws.Range("E40").Value = 7

# Row 41
$ws.Range("E41").Value = 19

# Row 48
$ws.Range("E48").Value = 13
$ws.Range("F48").Value = 7
$ws.Range("H48").Value = 7

# Row 49
$ws.Range("E49").Value = 38
$ws.Range("F49").Value = 18
$ws.Range("H49").Value = 18

# Row 61
$ws.Range("E61").Value = 17

# Row 62
$ws.Range("E62").Value = 18

# Row 68
$ws.Range("E68").Value = 9
$ws.Range("F68").Value = 3
$ws.Range("H68").Value = 3

# Row 78
$ws.Range("E78").Value = 14

# Row 83
$ws.Range("E83").Value = 6
